# View Html Report - View Stats - View Test Run Stats - Area Stats No Areas
#
# The "Steps" sheet has an "Area Stats" block (rows 93-99) that previously
# only covered a single pseudo-area ("Areas 1"). This edit renames that
# area to "Area 3" and appends a second area block ("Feature 7") with its
# own Scenarios/Steps counts in rows 101-105. It also switches the
# workbook to manual calculation and updates the sheet's view/selection
# state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Steps")
$ws.Activate()

# Rename the existing area header from "Areas 1" to "Area 3".
$ws.Range("A93").Value = "Area 3"

# New "Feature 7" area-stats block.
$ws.Range("A101").Value = "Feature 7"
$ws.Range("B101").Value = "Passing"
$ws.Range("C101").Value = "Skipped"
$ws.Range("D101").Value = "Failed"
$ws.Range("E101").Value = "Total"

$ws.Range("A102").Value = "Scenarios"
$ws.Range("B102").Value = 1
$ws.Range("C102").Value = 1
$ws.Range("D102").Value = 1
$ws.Range("E102").Value = 3

$ws.Range("B103").Formula = "=B102/E102"
$ws.Range("C103").Formula = "=C102/E102"
$ws.Range("D103").Formula = "=D102/E102"
$ws.Range("B103:D103").NumberFormat = "0.0000000000000%"

$ws.Range("A104").Value = "Steps"
$ws.Range("B104").Value = 4
$ws.Range("C104").Value = 4
$ws.Range("D104").Value = 1
$ws.Range("E104").Value = 9

$ws.Range("B105").Formula = "=B104/E104"
$ws.Range("C105").Formula = "=C104/E104"
$ws.Range("D105").Formula = "=D104/E104"
$ws.Range("B105:D105").NumberFormat = "0.0000000000000%"

# Scroll the frozen pane back to the top and move the selection.
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("G58").Select()

# Switch the workbook to manual calculation.
$excel.Calculation = -4135
